$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2021-Q1" sheet (currently the 2nd sheet) so that
#    we end up with three sheets: 总计, (old 2021-Q1 data), (copy, new 2021-Q1).
#    We will turn the original (2nd position) into the new "2022-Q4" sheet and
#    keep the copy (3rd position) as the untouched "2021-Q1" sheet.
# ---------------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(2)
$q1Sheet.Copy($null, $q1Sheet)

$newQ4 = $wb.Worksheets.Item(2)
$newQ1 = $wb.Worksheets.Item(3)

$newQ4.Name = "2022-Q4"
$newQ1.Name = "2021-Q1"

# ---------------------------------------------------------------------------
# 2. Populate the "2022-Q4" sheet with the new fund-holding data, overwriting
#    what used to be the "2021-Q1" numbers.
# ---------------------------------------------------------------------------

# Headers (row 1)
$newQ4.Range("B1").Value = "基金代码"
$newQ4.Range("C1").Value = "基金名称"
$newQ4.Range("D1").Value = "基金规模"
$newQ4.Range("E1").Value = "股票总仓位"
$newQ4.Range("F1").Value = "仓位占比"
$newQ4.Range("G1").Value = "持有市值(亿元)"
$newQ4.Range("H1").Value = "仓位排名"

# Row 2
$newQ4.Range("A2").Value = 0
$newQ4.Range("B2").NumberFormat = "@"
$newQ4.Range("B2").Value = "013680"
$newQ4.Range("C2").Value = "华安品质甄选混合A"
$newQ4.Range("D2").NumberFormat = "@"
$newQ4.Range("D2").Value = "12.95"
$newQ4.Range("E2").NumberFormat = "@"
$newQ4.Range("E2").Value = "73.22"
$newQ4.Range("F2").NumberFormat = "@"
$newQ4.Range("F2").Value = "1.93"
$newQ4.Range("G2").NumberFormat = "@"
$newQ4.Range("G2").Value = "0.2499"
$newQ4.Range("H2").Value = 4

# Row 3
$newQ4.Range("A3").Value = 1
$newQ4.Range("B3").NumberFormat = "@"
$newQ4.Range("B3").Value = "013681"
$newQ4.Range("C3").Value = "华安品质甄选混合C"
$newQ4.Range("D3").NumberFormat = "@"
$newQ4.Range("D3").Value = "5.10"
$newQ4.Range("E3").NumberFormat = "@"
$newQ4.Range("E3").Value = "73.22"
$newQ4.Range("F3").NumberFormat = "@"
$newQ4.Range("F3").Value = "1.93"
$newQ4.Range("G3").NumberFormat = "@"
$newQ4.Range("G3").Value = "0.0984"
$newQ4.Range("H3").Value = 4

# ---------------------------------------------------------------------------
# 3. Match the formatting used on the summary ("总计") sheet: header row /
#    first column use the same cell style, and the page margins line up with
#    sheet 1 as well.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Range("B1").Copy()
$newQ4.Range("B1:H1").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$newQ4.Range("A2:A3").PasteSpecial(-4122)

$newQ4.PageSetup.LeftMargin = 54
$newQ4.PageSetup.RightMargin = 54
$newQ4.PageSetup.TopMargin = 72
$newQ4.PageSetup.BottomMargin = 72
$newQ4.PageSetup.HeaderMargin = 36
$newQ4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 4. Update the "总计" summary sheet: row 2 now reports the new 2022-Q4
#    figures, and the old 2021-Q1 figures move down to a new row 3.
# ---------------------------------------------------------------------------
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q1"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.12

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.35
